$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.208684666666667
$ws.Range("H2").Value = 3.626054
$ws.Range("I2").Value = 0.01462795763842055
$ws.Range("J2").Value = 0.01462795763842055
$ws.Range("M2").Value = 4.063699000000001
$ws.Range("N2").Value = 12.191097
$ws.Range("Q2").Value = 4.911730671248668
$ws.Range("R2").Value = 44.205576041238
$ws.Range("S2").Value = 0.01462795763842055
$ws.Range("T2").Value = 0.01462795763842055

# Row 3
$ws.Range("I3").Value = 0.0626664797952065
$ws.Range("J3").Value = 0.06266647979520648
$ws.Range("M3").Value = 4.063699000000001
$ws.Range("N3").Value = 12.191097
$ws.Range("Q3").Value = 21.04195804210267
$ws.Range("R3").Value = 189.377622378924
$ws.Range("S3").Value = 0.0626664797952065
$ws.Range("T3").Value = 0.06266647979520648

# Row 4
$ws.Range("G4").Value = 76.16218566666667
$ws.Range("H4").Value = 228.486557
$ws.Range("I4").Value = 0.9217434921665711
$ws.Range("J4").Value = 0.921743492166571
$ws.Range("M4").Value = 4.063699000000001
$ws.Range("N4").Value = 12.191097
$ws.Range("Q4").Value = 309.5001977314477
$ws.Range("R4").Value = 2785.501779583029
$ws.Range("S4").Value = 0.9217434921665711
$ws.Range("T4").Value = 0.921743492166571

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07949433333333333
$ws.Range("H5").Value = 0.238483
$ws.Range("I5").Value = 0.0009620703998019471
$ws.Range("J5").Value = 0.000962070399801947
$ws.Range("M5").Value = 4.063699000000001
$ws.Range("N5").Value = 12.191097
$ws.Range("Q5").Value = 0.3230410428723334
$ws.Range("R5").Value = 2.907369385851
$ws.Range("S5").Value = 0.0009620703998019471
$ws.Range("T5").Value = 0.000962070399801947
